$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 9.040999999999999
$ws.Range("B7").Value = 5.230800000000003
$ws.Range("D7").Value = -7.387699999999995
$ws.Range("D15").Value = -8.508899999999995
$ws.Range("B16").Value = 7.098399999999997
$ws.Range("D21").Value = -8.531499999999999
$ws.Range("D22").Value = -7.981000000000006
$ws.Range("D23").Value = -7.178099999999997
$ws.Range("B28").Value = 6.0374
$ws.Range("B29").Value = 5.025700000000001
$ws.Range("B32").Value = 7.578899999999993
$ws.Range("D34").Value = -7.934300000000004
$ws.Range("B40").Value = 8.638000000000002
$ws.Range("D43").Value = -8.281300000000002
$ws.Range("D45").Value = -7.9138
$ws.Range("D50").Value = -8.172199999999997
$ws.Range("D51").Value = -7.424299999999999
$ws.Range("B52").Value = 5.105200000000001
$ws.Range("B57").Value = 4.964299999999995
$ws.Range("B66").Value = 5.555399999999999
$ws.Range("D66").Value = -7.374899999999999
$ws.Range("D67").Value = -6.351900000000001
$ws.Range("D79").Value = -6.273700000000002
$ws.Range("D84").Value = -8.910700000000002
$ws.Range("D92").Value = -6.450800000000002
$ws.Range("D97").Value = -8.412300000000002
$ws.Range("B100").Value = 5.894199999999999
